# Versuchsplan - Exp. 8 Drahtermuedung
# "Drahtsorten in rot & silber im Versuchsprotokoll geaendert"
#
# The "Drahtsorte" column (D) previously distinguished wires by material
# ("Eisen" / "Stahl"). Replace those values with the new colour-coded
# designations ("rot" / "silber") used throughout the protocol.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # column D = "Drahtsorte"
    $val = $cell.Value2

    if ($val -eq "Eisen") {
        $cell.Value2 = "rot"
    } elseif ($val -eq "Stahl") {
        $cell.Value2 = "silber"
    }
}

# Cosmetic touch-up: row 8's "Tag" label cell (B8) was missing the thin
# right border that every other row in that column carries - restore it
# so the table grid is consistent again.
$ws.Range("B8").Borders.Item(10).LineStyle = 1
$ws.Range("B8").Borders.Item(10).Weight = 2

# Move the active selection to C13 (also resets the saved scroll/view
# position back to the top of the sheet).
[void]$ws.Range("C13").Select()
